$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dev Log")
$ws.Rows(4).Insert()
